# Updated cryptos list values (Coin/Link/Price/Volume(1h)) per the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All of these columns (B: Coin, C: Link,
# D: Price, E: Volume(1h)) are stored as text in the sheet, including price strings
# that look numeric (e.g. "29.299.76", "0.9997"), so we force the cell format to
# Text ("@") before assigning -- otherwise Excel would coerce numeric-looking
# strings into actual numbers and drop things like trailing zeros.
$updates = [ordered]@{
    'D2' = '29.299.76'
    'E2' = '  -0.92%  '
    'D3' = '1.838.75'
    'E3' = '  -0.64%  '
    'D4' = '0.9997'
    'E4' = '  +0.08%  '
    'D5' = '238.94'
    'E5' = '  -0.69%  '
    'D6' = '0.6248'
    'E6' = '  -0.87%  '
    'E7' = '  +0.08%  '
    'D8' = '0.07359'
    'E8' = '  -1.45%  '
    'D9' = '0.2888'
    'E9' = '  -0.87%  '
    'D10' = '24.79'
    'E10' = '  -1.19%  '
    'D11' = '0.07713'
    'E11' = '  -0.44%  '
    'D12' = '1.840.56'
    'E12' = '  -0.65%  '
    'D13' = '4.947'
    'E13' = '  -1.53%  '
    'B14' = 'Polygon'
    'C14' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D14' = '0.6626'
    'E14' = '  -3.09%  '
    'B15' = 'ShibaInu'
    'C15' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    'D15' = '0.00001046'
    'E15' = '  +1.82%  '
    'E16' = '  -1.67%  '
    'D17' = '6.243'
    'E17' = '  -1.41%  '
    'D18' = '29.323.04'
    'E18' = '  -0.89%  '
    'D19' = '235.59'
    'E19' = '  +2.28%  '
    'E20' = '  -1.45%  '
    'E21' = '  +0.08%  '
    'D22' = '7.257'
    'E22' = '  -3.55%  '
    'E23' = '  +0.07%  '
    'D24' = '157.39'
    'E24' = '  -1.22%  '
    'D25' = '8.410'
    'E25' = '  -1.31%  '
    'D26' = '0.1335'
    'E26' = '  -2.37%  '
    'D27' = '17.25'
    'E27' = '  -1.84%  '
    'D28' = '0.07142'
    'E28' = '  +8.46%  '
    'D29' = '1.479'
    'E29' = '  +0.36%  '
    'E30' = '  -0.63%  '
    'B31' = 'Filecoin'
    'C31' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D31' = '4.017'
    'E31' = '  -2.26%  '
    'B32' = 'InternetComputer(DFINITY)'
    'C32' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D32' = '4.020'
    'E32' = '  -1.86%  '
    'D33' = '1.150'
    'E33' = '  +0.67%  '
    'D34' = '1.789'
    'E34' = '  -3.33%  '
    'D35' = '0.6898'
    'E35' = '  -1.47%  '
    'D36' = '2.582'
    'E36' = '  +0.57%  '
    'D37' = '0.01824'
    'E37' = '  -2.45%  '
    'B38' = 'MXToken'
    'C38' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D38' = '2.784'
    'E38' = '  -1.96%  '
    'B39' = 'Maker'
    'C39' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'D39' = '1.233.49'
    'E39' = '  -1.72%  '
    'D40' = '6.739'
    'E40' = '  -0.79%  '
    'D41' = '0.9440'
    'E41' = '  +0.80%  '
    'D42' = '1.002'
    'E42' = '  +0.04%  '
    'D43' = '1.996.81'
    'E43' = '  -0.25%  '
    'D44' = '101.24'
    'E44' = '  -0.17%  '
    'D45' = '65.05'
    'E45' = '  -1.85%  '
    'D46' = '0.00000000118'
    'E46' = '  +5.53%  '
    'D47' = '6.915'
    'E47' = '  -2.80%  '
    'D48' = '1.681'
    'E48' = '  -2.94%  '
    'D49' = '8.878'
    'E49' = '  -1.57%  '
    'D50' = '0.1127'
    'E50' = '  -2.66%  '
    'D51' = '0.3875'
    'E51' = '  -1.69%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
